# Slide 2 ("Актуальність теми"), body placeholder: trim the parenthetical
# list down to "(финансовые графики, количество игр)" — the trailing
# ", улучшение ши интелектов" is removed and the run that used to hold
# ", " before it is retyped as the closing ")" in English.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(2)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# Characters() uses 1-based indices, same as classic PowerPoint VBA.
# The run ", " (chars 199-200) together with the following runs
# "улучшение" " " "ши" " " "интелектов" (chars 201-223) collapse into a
# single ")" character that closes the sentence.
$target = $tr.Characters(199, 25)
$target.Text = ")"
$target.Font.LanguageID = "en-US"
